# Auto-generated Excel COM-interop script applying the Titan_Profits price-update diff
# Recomputes H/I/J/K/L/M/N columns for specific (sheet,row) pairs as captured in the diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 227.15384
$ws.Range("I5").Value = 51.57143
$ws.Range("K5").Value = 51.57143
$ws.Range("M5").Value = 63.42857

$ws.Range("H32").Value = 1683.6666
$ws.Range("I32").Value = 750
$ws.Range("J32").Value = 2150.5
$ws.Range("K32").Value = 750
$ws.Range("L32").Value = 2150.5
$ws.Range("M32").Value = -424
$ws.Range("N32").Value = -2802.5

$ws.Range("H76").Value = 3834345.2
$ws.Range("I76").Value = 4118093
$ws.Range("K76").Value = 4118093
$ws.Range("M76").Value = -4117778

$ws.Range("H79").Value = 3834345.2
$ws.Range("I79").Value = 4118093
$ws.Range("K79").Value = 4118093
$ws.Range("M79").Value = -4117001

$ws.Range("H105").Value = 500671
$ws.Range("J105").Value = 500671
$ws.Range("L105").Value = 500671
$ws.Range("N105").Value = -507659

$ws.Range("H132").Value = 218276.7
$ws.Range("I132").Value = 234267.27
$ws.Range("K132").Value = 702801.8099999999
$ws.Range("M132").Value = -700271.8099999999

$ws.Range("H135").Value = 2202.6843
$ws.Range("I135").Value = 2161
$ws.Range("J135").Value = 2359
$ws.Range("K135").Value = 19449
$ws.Range("L135").Value = 21231
$ws.Range("M135").Value = -16914
$ws.Range("N135").Value = -26301

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1402.1522
$ws.Range("I61").Value = 938.3946999999999
$ws.Range("J61").Value = 3605
$ws.Range("K61").Value = 938.3946999999999
$ws.Range("L61").Value = 3605
$ws.Range("M61").Value = -726.3946999999999
$ws.Range("N61").Value = -4029

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H122").Value = 1934.8235
$ws.Range("I122").Value = 1663.4736
$ws.Range("J122").Value = 2278.5334
$ws.Range("K122").Value = 4990.4208
$ws.Range("L122").Value = 6835.600199999999
$ws.Range("M122").Value = -2540.4208
$ws.Range("N122").Value = -11735.6002

$ws.Range("H136").Value = 1402.1522
$ws.Range("I136").Value = 938.3946999999999
$ws.Range("J136").Value = 3605
$ws.Range("K136").Value = 2815.1841
$ws.Range("L136").Value = 10815
$ws.Range("M136").Value = -265.1840999999999
$ws.Range("N136").Value = -15915

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3890.3333
$ws.Range("I132").Value = 3208.8
$ws.Range("J132").Value = 4377.143
$ws.Range("K132").Value = 9626.400000000001
$ws.Range("L132").Value = 13131.429
$ws.Range("M132").Value = -7096.400000000001
$ws.Range("N132").Value = -18191.429

$ws.Range("H133").Value = 17574.5
$ws.Range("I133").Value = 20296
$ws.Range("J133").Value = 17185.715
$ws.Range("K133").Value = 20296
$ws.Range("L133").Value = 17185.715
$ws.Range("M133").Value = -17766
$ws.Range("N133").Value = -22245.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1376.3103
$ws.Range("I5").Value = 550.0833
$ws.Range("J5").Value = 1959.5294
$ws.Range("K5").Value = 1650.2499
$ws.Range("L5").Value = 5878.5882
$ws.Range("M5").Value = -1538.2499
$ws.Range("N5").Value = -6102.5882

$ws.Range("H18").Value = 1829
$ws.Range("I18").Value = 217.5
$ws.Range("J18").Value = 3977.6667
$ws.Range("K18").Value = 652.5
$ws.Range("L18").Value = 11933.0001
$ws.Range("M18").Value = -483.5
$ws.Range("N18").Value = -12271.0001

$ws.Range("H70").Value = 5000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 15000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -15630

$ws.Range("H73").Value = 5000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 15000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -17184

$ws.Range("H75").Value = 2554.5
$ws.Range("I75").Value = 2171
$ws.Range("J75").Value = 3705
$ws.Range("K75").Value = 6513
$ws.Range("L75").Value = 11115
$ws.Range("M75").Value = -5515
$ws.Range("N75").Value = -13111

$ws.Range("H78").Value = 2554.5
$ws.Range("I78").Value = 2171
$ws.Range("J78").Value = 3705
$ws.Range("K78").Value = 19539
$ws.Range("L78").Value = 33345
$ws.Range("M78").Value = -14547
$ws.Range("N78").Value = -43329

$ws.Range("H103").Value = 1793.2106
$ws.Range("I103").Value = 1262.1428
$ws.Range("J103").Value = 2103
$ws.Range("K103").Value = 3786.4284
$ws.Range("L103").Value = 6309
$ws.Range("M103").Value = -2907.4284
$ws.Range("N103").Value = -8067

$ws.Range("H121").Value = 614.5454999999999
$ws.Range("I121").Value = 120
$ws.Range("K121").Value = 360
$ws.Range("M121").Value = 950

$ws.Range("H135").Value = 1376.3103
$ws.Range("I135").Value = 550.0833
$ws.Range("J135").Value = 1959.5294
$ws.Range("K135").Value = 4950.7497
$ws.Range("L135").Value = 17635.7646
$ws.Range("M135").Value = -2415.7497
$ws.Range("N135").Value = -22705.7646

$ws.Range("H139").Value = 2104.7568
$ws.Range("I139").Value = 1888.875
$ws.Range("J139").Value = 3486.4
$ws.Range("K139").Value = 5666.625
$ws.Range("L139").Value = 10459.2
$ws.Range("M139").Value = -526.625
$ws.Range("N139").Value = -20739.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5839.147
$ws.Range("I70").Value = 5912.731
$ws.Range("J70").Value = 5600
$ws.Range("K70").Value = 5912.731
$ws.Range("L70").Value = 5600
$ws.Range("M70").Value = -5642.731
$ws.Range("N70").Value = -6140

$ws.Range("H73").Value = 5839.147
$ws.Range("I73").Value = 5912.731
$ws.Range("J73").Value = 5600
$ws.Range("K73").Value = 5912.731
$ws.Range("L73").Value = 5600
$ws.Range("M73").Value = -4976.731
$ws.Range("N73").Value = -7472

$ws.Range("H93").Value = 27849.5
$ws.Range("J93").Value = 27849.5
$ws.Range("L93").Value = 27849.5
$ws.Range("N93").Value = -31593.5

$ws.Range("H97").Value = 775
$ws.Range("I97").Value = 609.5833
$ws.Range("J97").Value = 1058.5714
$ws.Range("K97").Value = 609.5833
$ws.Range("L97").Value = 1058.5714
$ws.Range("M97").Value = -113.5833
$ws.Range("N97").Value = -2050.5714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 6250
$ws.Range("I5").Value = 7500
$ws.Range("J5").Value = 5000
$ws.Range("K5").Value = 7500
$ws.Range("L5").Value = 5000
$ws.Range("M5").Value = -7387
$ws.Range("N5").Value = -5226

$ws.Range("H22").Value = 10550.091
$ws.Range("I22").Value = 701
$ws.Range("J22").Value = 11535
$ws.Range("K22").Value = 701
$ws.Range("L22").Value = 11535
$ws.Range("M22").Value = -406
$ws.Range("N22").Value = -12125

$ws.Range("H27").Value = 10550.091
$ws.Range("I27").Value = 701
$ws.Range("J27").Value = 11535
$ws.Range("K27").Value = 701
$ws.Range("L27").Value = 11535
$ws.Range("M27").Value = -594
$ws.Range("N27").Value = -11749

$ws.Range("H31").Value = 975.7143
$ws.Range("I31").Value = 855
$ws.Range("K31").Value = 855
$ws.Range("M31").Value = -607

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 24647.812
$ws.Range("I14").Value = 50000
$ws.Range("K14").Value = 50000
$ws.Range("M14").Value = -49832

$ws.Range("H23").Value = 627.5
$ws.Range("I23").Value = 627.5
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 627.5
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -398.5
$ws.Range("N23").ClearContents()

$ws.Range("H132").Value = 13515726
$ws.Range("I132").Value = 16130714
$ws.Range("J132").Value = 4954.5
$ws.Range("K132").Value = 48392142
$ws.Range("L132").Value = 14863.5
$ws.Range("M132").Value = -48389612
$ws.Range("N132").Value = -19923.5

$ws.Range("H136").Value = 6824283
$ws.Range("I136").Value = 7599239
$ws.Range("J136").Value = 4668.4
$ws.Range("K136").Value = 22797717
$ws.Range("L136").Value = 14005.2
$ws.Range("M136").Value = -22795167
$ws.Range("N136").Value = -19105.2
